# Update Sheets via scheduled runner: refresh market-data columns (H-N)
# for the specific Leve rows identified in the upstream data pull.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1699.25
$ws.Range("I32").Value = 1400
$ws.Range("J32").Value = 1998.5
$ws.Range("K32").Value = 1400
$ws.Range("L32").Value = 1998.5
$ws.Range("M32").Value = -1074
$ws.Range("N32").Value = -2650.5
$ws.Range("H53").Value = 585.8125
$ws.Range("I53").Value = 339.4
$ws.Range("J53").Value = 996.5
$ws.Range("K53").Value = 339.4
$ws.Range("L53").Value = 996.5
$ws.Range("M53").Value = 297.6
$ws.Range("N53").Value = -2270.5
$ws.Range("H88").Value = 2450.6428
$ws.Range("I88").Value = 641.5
$ws.Range("J88").Value = 2752.1667
$ws.Range("K88").Value = 641.5
$ws.Range("L88").Value = 2752.1667
$ws.Range("M88").Value = -235.5
$ws.Range("N88").Value = -3564.1667
$ws.Range("H91").Value = 2450.6428
$ws.Range("I91").Value = 641.5
$ws.Range("J91").Value = 2752.1667
$ws.Range("K91").Value = 641.5
$ws.Range("L91").Value = 2752.1667
$ws.Range("M91").Value = 762.5
$ws.Range("N91").Value = -5560.1667
$ws.Range("H107").Value = 221.875
$ws.Range("I107").Value = 230.16667
$ws.Range("K107").Value = 230.16667
$ws.Range("M107").Value = 1689.83333
$ws.Range("H132").Value = 4349.5713
$ws.Range("I132").Value = 2128.1428
$ws.Range("K132").Value = 6384.428400000001
$ws.Range("M132").Value = -3854.428400000001

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4208
$ws.Range("I45").Value = 2664.5
$ws.Range("J45").Value = 5442.8
$ws.Range("K45").Value = 2664.5
$ws.Range("L45").Value = 5442.8
$ws.Range("M45").Value = -2287.5
$ws.Range("N45").Value = -6196.8
$ws.Range("H46").Value = 18987.5
$ws.Range("J46").Value = 18983.334
$ws.Range("L46").Value = 18983.334
$ws.Range("N46").Value = -19621.334
$ws.Range("H61").Value = 4289.2
$ws.Range("I61").Value = 4154.778
$ws.Range("J61").Value = 5499
$ws.Range("K61").Value = 4154.778
$ws.Range("L61").Value = 5499
$ws.Range("M61").Value = -3942.778
$ws.Range("N61").Value = -5923
$ws.Range("H74").Value = 754.875
$ws.Range("I74").Value = 754.875
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 754.875
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 119.125
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 754.875
$ws.Range("I77").Value = 754.875
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3774.375
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 593.625
$ws.Range("N77").ClearContents()
$ws.Range("H96").Value = 16172
$ws.Range("J96").Value = 16172
$ws.Range("L96").Value = 16172
$ws.Range("N96").Value = -21664
$ws.Range("H97").Value = 7645
$ws.Range("J97").Value = 15094.5
$ws.Range("L97").Value = 15094.5
$ws.Range("N97").Value = -16086.5
$ws.Range("H132").Value = 1488
$ws.Range("I132").Value = 1488
$ws.Range("K132").Value = 4464
$ws.Range("M132").Value = -1934
$ws.Range("H136").Value = 4289.2
$ws.Range("I136").Value = 4154.778
$ws.Range("J136").Value = 5499
$ws.Range("K136").Value = 12464.334
$ws.Range("L136").Value = 16497
$ws.Range("M136").Value = -9914.334000000001
$ws.Range("N136").Value = -21597

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2085.5454
$ws.Range("I86").Value = 1877.4286
$ws.Range("K86").Value = 1877.4286
$ws.Range("M86").Value = -754.4286
$ws.Range("H89").Value = 2085.5454
$ws.Range("I89").Value = 1877.4286
$ws.Range("K89").Value = 9387.143
$ws.Range("M89").Value = -3771.143
$ws.Range("H99").Value = 2682.7778
$ws.Range("I99").Value = 2424
$ws.Range("K99").Value = 2424
$ws.Range("M99").Value = -926
$ws.Range("H134").Value = 6000
$ws.Range("J134").Value = 6000
$ws.Range("L134").Value = 18000
$ws.Range("N134").Value = -23070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2200
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 2200
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 5216.1665
$ws.Range("I132").Value = 5405.353
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 16216.059
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -13686.059
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 2060.2856
$ws.Range("J134").Value = 2499.5
$ws.Range("L134").Value = 7498.5
$ws.Range("N134").Value = -12568.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 113.2
$ws.Range("I2").Value = 97.5
$ws.Range("K2").Value = 585
$ws.Range("M2").Value = -472
$ws.Range("H107").Value = 859.6667
$ws.Range("J107").Value = 1038
$ws.Range("L107").Value = 3114
$ws.Range("N107").Value = -6954

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 529.52
$ws.Range("I97").Value = 505.42856
$ws.Range("J97").Value = 656
$ws.Range("K97").Value = 505.42856
$ws.Range("L97").Value = 656
$ws.Range("M97").Value = -9.428560000000004
$ws.Range("N97").Value = -1648
$ws.Range("H122").Value = 2142.7144
$ws.Range("I122").Value = 1999.8334
$ws.Range("K122").Value = 5999.5002
$ws.Range("M122").Value = -3549.5002
$ws.Range("H126").Value = 8669.666999999999
$ws.Range("I126").Value = 8005
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 24015
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -21545
$ws.Range("N126").Value = -34937
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2190.1538
$ws.Range("I40").Value = 1789.3334
$ws.Range("K40").Value = 1789.3334
$ws.Range("M40").Value = -1653.3334
$ws.Range("H46").Value = 1574
$ws.Range("J46").Value = 1642.7858
$ws.Range("L46").Value = 1642.7858
$ws.Range("N46").Value = -2018.7858
$ws.Range("H55").Value = 661.7143
$ws.Range("I55").Value = 175.33333
$ws.Range("J55").Value = 856.26666
$ws.Range("K55").Value = 175.33333
$ws.Range("L55").Value = 856.26666
$ws.Range("M55").Value = -2.333329999999989
$ws.Range("N55").Value = -1202.26666
$ws.Range("H100").Value = 2128.2856
$ws.Range("I100").Value = 2128.2856
$ws.Range("K100").Value = 2128.2856
$ws.Range("M100").Value = -1587.2856

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 58849.668
$ws.Range("I132").Value = 43477.64
$ws.Range("J132").Value = 251000
$ws.Range("K132").Value = 130432.92
$ws.Range("L132").Value = 753000
$ws.Range("M132").Value = -127902.92
$ws.Range("N132").Value = -758060
$ws.Range("H136").Value = 4495
$ws.Range("I136").Value = 4495
$ws.Range("K136").Value = 13485
$ws.Range("M136").Value = -10935

